$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the connector diameter (cn_D) equation in B7:
# "13 mm + 1.5*tol" -> "13 mm + 2*tol"
$ws.Range("B7").Value = "13 mm + 2*tol"

# Reflect the final selection recorded in the workbook (B21)
$ws.Range("B21").Select()
